# Edit slide 1, shape "TextBox 18" text:
#   1) "W {Item Width (inch)} cm x D ..."  ->  "W {Item Width(inch)} cm x D ..."
#      (drop the space between "Width" and "(inch)")
#   2) "Suggested Retail AU: ${Retail AUD}" -> "Suggested Retail AU: ${Retail AU}"
#      (drop the "D" in "AUD}")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 18")
$tr = $sh.TextFrame.TextRange

# --- Fix #1: "Width (inch)..." -> "Width(inch)..." -------------------------
$full = $tr.Text
$oldSeg1 = "Width (inch)} cm x D {Item Depth (inch)} cm x H {Item Height (inch)} cm "
$newSeg1 = "Width(inch)} cm x D {Item Depth (inch)} cm x H {Item Height (inch)} cm "
$idx1 = $full.IndexOf($oldSeg1)
if ($idx1 -ge 0) {
    $seg1 = $tr.Characters($idx1 + 1, $oldSeg1.Length)
    $seg1.Text = $newSeg1
}

# --- Fix #2: "AUD}" -> "AU}" ------------------------------------------------
$full = $tr.Text
$oldSeg2 = "AUD}"
$newSeg2 = "AU}"
$idx2 = $full.IndexOf($oldSeg2)
if ($idx2 -ge 0) {
    $seg2 = $tr.Characters($idx2 + 1, $oldSeg2.Length)
    $seg2.Text = $newSeg2
}
